$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data (empadronador, total_registros) sorted descending by total_registros,
# reflecting the new counts from the latest upload.
$data = @(
    @("TOLENTINO VASQUEZ DIANA KATHERYN", 60),
    @("JULCA VALENZUELA CINTIA KARYN", 58),
    @("SANCHEZ CORTEZ LEYLA DIANA", 54),
    @("CARRILLO MARTÍNEZ HEIDY NAYELI", 51),
    @("VALER VEGA PATRICIA GERALDINE", 51),
    @("DE LA CRUZ BENITES RICHARD ALEXANDER", 50),
    @("YZQUIERDO CARHUATANTA LEYDY YANELA", 48),
    @("RODRIGUEZ RUBIO SANDRA MABEL", 48),
    @("PONCE VILLANUEVA CARMEN ISABEL", 46),
    @("ARENAS ZAVALA ANDYELA PATRICIA ISIDORA", 45),
    @("REYES RODRIGUEZ JEISSON STEVEN", 44),
    @("GASLAC GUTIERREZ FRANK JHORDY", 43),
    @("RUBIO MARIÑOS GISELA JUDITH", 42),
    @("PIERINA NAGIELLY SANDOVAL CONTRERAS", 41),
    @("CYNTHIA RODRIGUEZ LECCA", 40),
    @("SEGURA ASTO YAMILET ANTONELA", 38),
    @("GUZMAN ZAVALETA CECILIA MARISOL", 33),
    @("RODRIGUEZ VASQUEZ WALTER", 30),
    @("LEON VERA MELISSA FIORELLA", 25),
    @("RUTH MELISSA RAMIREZ VELEZMORO", 16)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
